$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.04904614488327752
$ws.Range("D2").Value = 0.1367748167803953
$ws.Range("E2").Value = 0.1464392491620963
$ws.Range("F2").Value = 2.05016538974661
$ws.Range("G2").Value = 0.00243455647104717
$ws.Range("J2").Value = 0.200971130535109
$ws.Range("N2").Value = 2.873581687696401
$ws.Range("O2").Value = 5.513759881882606
$ws.Range("C3").Value = 0.0435296659521498
$ws.Range("D3").Value = 0.1371026173879173
$ws.Range("E3").Value = 0.1446177847971448
$ws.Range("F3").Value = 1.980028762343835
$ws.Range("G3").Value = 0.002440342032222053
$ws.Range("J3").Value = 0.1953332899244486
$ws.Range("N3").Value = 2.562605684679681
$ws.Range("O3").Value = 5.291219219213872
$ws.Range("C4").Value = 0.04016292325233906
$ws.Range("D4").Value = 0.1373433076148487
$ws.Range("E4").Value = 0.1435801722260521
$ws.Range("F4").Value = 1.938296182754442
$ws.Range("G4").Value = 0.002444081947155112
$ws.Range("J4").Value = 0.1920013065015809
$ws.Range("N4").Value = 2.371325805375818
$ws.Range("O4").Value = 5.158232048430705
$ws.Range("C5").Value = 0.03879591925561954
$ws.Range("D5").Value = 0.1374513150846894
$ws.Range("E5").Value = 0.1431775837974101
$ws.Range("F5").Value = 1.921622456926414
$ws.Range("G5").Value = 0.002445653326600058
$ws.Range("J5").Value = 0.190675899207605
$ws.Range("N5").Value = 2.293303068607429
$ws.Range("O5").Value = 5.104950414272651
$ws.Range("C6").Value = 0.0385692249258085
$ws.Range("D6").Value = 0.1374698494004569
$ws.Range("E6").Value = 0.143111955325697
$ws.Range("F6").Value = 1.918873824376234
$ws.Range("G6").Value = 0.00244591711691869
$ws.Range("J6").Value = 0.1904577674833092
$ws.Range("N6").Value = 2.280343261403573
$ws.Range("O6").Value = 5.096157901393042
$ws.Range("C7").Value = 0.0401444674392053
$ws.Range("D7").Value = 0.1373447240420163
$ws.Range("E7").Value = 0.1435746608810646
$ws.Range("F7").Value = 1.938069971089746
$ws.Range("G7").Value = 0.002444102947451032
$ws.Range("J7").Value = 0.1919833006640488
$ws.Range("N7").Value = 2.370273851392596
$ws.Range("O7").Value = 5.157509791371865
$ws.Range("C8").Value = 0.04713973221768697
$ws.Range("D8").Value = 0.1368796679006117
$ws.Range("E8").Value = 0.1457944037207461
$ws.Range("F8").Value = 2.02570454850482
$ws.Range("G8").Value = 0.002436512508591203
$ws.Range("J8").Value = 0.1990001632843672
$ws.Range("N8").Value = 2.766433886209825
$ws.Range("O8").Value = 5.436265195847056
$ws.Range("C9").Value = 0.06102711453344511
$ws.Range("D9").Value = 0.1362800148173093
$ws.Range("E9").Value = 0.1507915675396667
$ws.Range("F9").Value = 2.208232907395399
$ws.Range("G9").Value = 0.002423108042852531
$ws.Range("J9").Value = 0.2137993797248896
$ws.Range("N9").Value = 3.540180268007646
$ws.Range("O9").Value = 6.012248861075818
$ws.Range("C10").Value = 0.07134550460941114
$ws.Range("D10").Value = 0.1360293435454452
$ws.Range("E10").Value = 0.1548610410017872
$ws.Range("F10").Value = 2.349018070605325
$ws.Range("G10").Value = 0.002414151305390133
$ws.Range("J10").Value = 0.2253216067525585
$ws.Range("N10").Value = 4.10623028343673
$ws.Range("O10").Value = 6.453859072126306
$ws.Range("C11").Value = 0.07606729817744906
$ws.Range("D11").Value = 0.1359564470757348
$ws.Range("E11").Value = 0.1567999944200196
$ws.Range("F11").Value = 2.414554617585026
$ws.Range("G11").Value = 0.002410267899654984
$ws.Range("J11").Value = 0.2307078846769457
$ws.Range("N11").Value = 4.363110593465422
$ws.Range("O11").Value = 6.658884712929989
$ws.Range("C12").Value = 0.07785952777632588
$ws.Range("D12").Value = 0.1359347498808319
$ws.Range("E12").Value = 0.1575469301407395
$ws.Range("F12").Value = 2.439589128552171
$ws.Range("G12").Value = 0.002408824648859085
$ws.Range("J12").Value = 0.2327686092154835
$ws.Range("N12").Value = 4.460285735714251
$ws.Range("O12").Value = 6.737126306206733
$ws.Range("C13").Value = 0.07747334993821653
$ws.Range("D13").Value = 0.1359391601632964
$ws.Range("E13").Value = 0.1573854981983374
$ws.Range("F13").Value = 2.434187794826073
$ws.Range("G13").Value = 0.00240913426668442
$ws.Range("J13").Value = 0.2323238548572135
$ws.Range("N13").Value = 4.439361943450422
$ws.Range("O13").Value = 6.720248635271673
$ws.Range("C14").Value = 0.07621466101781493
$ws.Range("D14").Value = 0.1359545436921366
$ws.Range("E14").Value = 0.1568611903521031
$ws.Range("F14").Value = 2.416609852933618
$ws.Range("G14").Value = 0.002410148615865744
$ws.Range("J14").Value = 0.2308769982902845
$ws.Range("N14").Value = 4.371107314139522
$ws.Range("O14").Value = 6.665309566532869
$ws.Range("C15").Value = 0.07544422848580723
$ws.Range("D15").Value = 0.1359647355959055
$ws.Range("E15").Value = 0.1565416924997649
$ws.Range("F15").Value = 2.405871225399551
$ws.Range("G15").Value = 0.002410773486131744
$ws.Range("J15").Value = 0.2299935064924625
$ws.Range("N15").Value = 4.329286057409945
$ws.Range("O15").Value = 6.63173655416108
$ws.Range("C16").Value = 0.07103750082856664
$ws.Range("D16").Value = 0.1360349345466503
$ws.Range("E16").Value = 0.1547360980360182
$ws.Range("F16").Value = 2.344765360487287
$ws.Range("G16").Value = 0.002414408925553005
$ws.Range("J16").Value = 0.2249725343124993
$ws.Range("N16").Value = 4.089429168003562
$ws.Range("O16").Value = 6.440544070253338
$ws.Range("C17").Value = 0.06834138805399448
$ws.Range("D17").Value = 0.136088529159629
$ws.Range("E17").Value = 0.153650950788645
$ws.Range("F17").Value = 2.307663033581065
$ws.Range("G17").Value = 0.002416687968493563
$ws.Range("J17").Value = 0.2219295774007861
$ws.Range("N17").Value = 3.94211849063862
$ws.Range("O17").Value = 6.32431854201144
$ws.Range("C18").Value = 0.06679327129147339
$ws.Range("D18").Value = 0.1361232276543021
$ws.Range("E18").Value = 0.1530350553665372
$ws.Range("F18").Value = 2.286463095253197
$ws.Range("G18").Value = 0.002418016806669959
$ws.Range("J18").Value = 0.2201929597447361
$ws.Range("N18").Value = 3.857331695637754
$ws.Range("O18").Value = 6.257857391018604
$ws.Range("C19").Value = 0.06626954953216568
$ws.Range("D19").Value = 0.1361356412928245
$ws.Range("E19").Value = 0.1528279387180866
$ws.Range("F19").Value = 2.279309185331812
$ws.Range("G19").Value = 0.002418469824216481
$ws.Range("J19").Value = 0.2196073011980104
$ws.Range("N19").Value = 3.828614786364199
$ws.Range("O19").Value = 6.235421298198503
$ws.Range("C20").Value = 0.06862812184805023
$ws.Range("D20").Value = 0.1360824232132813
$ws.Range("E20").Value = 0.1537656118770911
$ws.Range("F20").Value = 2.31159809642412
$ws.Range("G20").Value = 0.002416443499604631
$ws.Range("J20").Value = 0.2222520949687379
$ws.Range("N20").Value = 3.95780600327754
$ws.Range("O20").Value = 6.336650659907093
$ws.Range("C21").Value = 0.07658425336667563
$ws.Range("D21").Value = 0.1359498649216775
$ws.Range("E21").Value = 0.1570148470531691
$ws.Range("F21").Value = 2.421767003689467
$ws.Range("G21").Value = 0.002409849937127243
$ws.Range("J21").Value = 0.2313014016787349
$ws.Range("N21").Value = 4.391158149571083
$ws.Range("O21").Value = 6.6814300756796
$ws.Range("C22").Value = 0.08180854445838293
$ws.Range("D22").Value = 0.1358976581555851
$ws.Range("E22").Value = 0.1592124517862601
$ws.Range("F22").Value = 2.495036253509966
$ws.Range("G22").Value = 0.00240569977638522
$ws.Range("J22").Value = 0.2373385043006238
$ws.Range("N22").Value = 4.67379181795809
$ws.Range("O22").Value = 6.910281055859514
$ws.Range("C23").Value = 0.07901794153349329
$ws.Range("D23").Value = 0.1359223743079596
$ws.Range("E23").Value = 0.1580327467680078
$ws.Range("F23").Value = 2.455814227775591
$ws.Range("G23").Value = 0.002407900290106358
$ws.Range("J23").Value = 0.2341050664453377
$ws.Range("N23").Value = 4.523002190001307
$ws.Range("O23").Value = 6.787814350591134
$ws.Range("C24").Value = 0.06849848355020072
$ws.Range("D24").Value = 0.1360851716060196
$ws.Range("E24").Value = 0.1537137487797153
$ws.Range("F24").Value = 2.309818647658062
$ws.Range("G24").Value = 0.002416553966016216
$ws.Range("J24").Value = 0.2221062448667368
$ws.Range("N24").Value = 3.950713976768498
$ws.Range("O24").Value = 6.331074193523932
$ws.Range("C25").Value = 0.05725082268372716
$ws.Range("D25").Value = 0.1364088544316537
$ws.Range("E25").Value = 0.1493701888532932
$ws.Range("F25").Value = 2.15769341987405
$ws.Range("G25").Value = 0.002426576956464862
$ws.Range("J25").Value = 0.2096829376187799
$ws.Range("N25").Value = 3.331249627311138
$ws.Range("O25").Value = 5.853232542874196
